$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 767, shifting the existing rows 767:808 down to 768:809.
$ws.Rows.Item(767).EntireRow.Insert()

# Populate the newly inserted row 767 with the new data point.
# Use a text number format while assigning the date-like value so Excel keeps it
# as plain text (matching the rest of the "date" column) instead of auto-converting
# it to a date serial number, then clear the format again so no stray style is left
# behind on the cell (the other cells in this column carry no explicit style either).
$ws.Range("A767").NumberFormat = "@"
$ws.Range("A767").Value = "2026/02/04"
$ws.Range("A767").ClearFormats()

$ws.Range("B767").Value = "水"
$ws.Range("C767").Value = 2
$ws.Range("D767").Value = 201
